$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.ClearFormats()
}

Set-TextValue $ws.Range('D2') '26.423.93'
Set-TextValue $ws.Range('E2') '  +1.28%  '
Set-TextValue $ws.Range('D3') '1.692.50'
Set-TextValue $ws.Range('E3') '  +1.42%  '
Set-TextValue $ws.Range('D4') '1.012'
Set-TextValue $ws.Range('E4') '  +0.90%  '
Set-TextValue $ws.Range('D5') '218.84'
Set-TextValue $ws.Range('E5') '  +1.23%  '
Set-TextValue $ws.Range('D6') '0.5484'
Set-TextValue $ws.Range('E6') '  +7.43%  '
Set-TextValue $ws.Range('D7') '1.012'
Set-TextValue $ws.Range('E7') '  +0.83%  '
Set-TextValue $ws.Range('D8') '0.2719'
Set-TextValue $ws.Range('E8') '  +1.12%  '
Set-TextValue $ws.Range('D9') '0.06474'
Set-TextValue $ws.Range('E9') '  +1.52%  '
Set-TextValue $ws.Range('D10') '22.09'
Set-TextValue $ws.Range('E10') '  +0.98%  '
Set-TextValue $ws.Range('D11') '0.07705'
Set-TextValue $ws.Range('E11') '  +3.61%  '
Set-TextValue $ws.Range('B12') 'Polkadot'
Set-TextValue $ws.Range('C12') 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
Set-TextValue $ws.Range('D12') '4.545'
Set-TextValue $ws.Range('E12') '  +0.67%  '
Set-TextValue $ws.Range('B13') 'WrappedEther'
Set-TextValue $ws.Range('C13') 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
Set-TextValue $ws.Range('D13') '1.679.82'
Set-TextValue $ws.Range('E13') '  +0.17%  '
Set-TextValue $ws.Range('D14') '0.5822'
Set-TextValue $ws.Range('E14') '  +0.15%  '
Set-TextValue $ws.Range('D15') '0.000008416'
Set-TextValue $ws.Range('E15') '  -0.85%  '
Set-TextValue $ws.Range('D16') '65.26'
Set-TextValue $ws.Range('E16') '  +1.74%  '
Set-TextValue $ws.Range('D17') '26.488.11'
Set-TextValue $ws.Range('E17') '  +2.26%  '
Set-TextValue $ws.Range('D18') '4.956'
Set-TextValue $ws.Range('E18') '  +0.61%  '
Set-TextValue $ws.Range('D19') '1.012'
Set-TextValue $ws.Range('E19') '  +0.87%  '
Set-TextValue $ws.Range('D20') '10.99'
Set-TextValue $ws.Range('E20') '  +1.64%  '
Set-TextValue $ws.Range('D21') '190.05'
Set-TextValue $ws.Range('E21') '  +0.19%  '
Set-TextValue $ws.Range('D22') '6.237'
Set-TextValue $ws.Range('E22') '  +0.72%  '
Set-TextValue $ws.Range('D23') '1.013'
Set-TextValue $ws.Range('E23') '  +0.87%  '
Set-TextValue $ws.Range('D24') '150.17'
Set-TextValue $ws.Range('E24') '  +3.85%  '
Set-TextValue $ws.Range('D25') '0.1303'
Set-TextValue $ws.Range('E25') '  +6.69%  '
Set-TextValue $ws.Range('D26') '7.900'
Set-TextValue $ws.Range('E26') '  +3.92%  '
Set-TextValue $ws.Range('D27') '15.73'
Set-TextValue $ws.Range('E27') '  +0.21%  '
Set-TextValue $ws.Range('D28') '1.422'
Set-TextValue $ws.Range('E28') '  +5.46%  '
Set-TextValue $ws.Range('D29') '0.06323'
Set-TextValue $ws.Range('E30') '  +1.02%  '
Set-TextValue $ws.Range('D31') '3.582'
Set-TextValue $ws.Range('E31') '  +0.39%  '
Set-TextValue $ws.Range('D32') '3.592'
Set-TextValue $ws.Range('E32') '  +1.99%  '
Set-TextValue $ws.Range('D33') '1.676'
Set-TextValue $ws.Range('E33') '  +0.64%  '
Set-TextValue $ws.Range('D34') '1.042'
Set-TextValue $ws.Range('E34') '  +2.45%  '
Set-TextValue $ws.Range('D35') '0.6213'
Set-TextValue $ws.Range('E35') '  +0.76%  '
Set-TextValue $ws.Range('D36') '2.414'
Set-TextValue $ws.Range('E36') '  +1.94%  '
Set-TextValue $ws.Range('E37') '  +1.51%  '
Set-TextValue $ws.Range('D38') '6.222'
Set-TextValue $ws.Range('E38') '  -0.62%  '
Set-TextValue $ws.Range('D39') '1.121.97'
Set-TextValue $ws.Range('E39') '  +2.32%  '
Set-TextValue $ws.Range('D40') '0.01640'
Set-TextValue $ws.Range('E40') '  +2.69%  '
Set-TextValue $ws.Range('D41') '0.8821'
Set-TextValue $ws.Range('E41') '  +1.89%  '
Set-TextValue $ws.Range('E42') '  +0.75%  '
Set-TextValue $ws.Range('D43') '100.99'
Set-TextValue $ws.Range('E43') '  -0.28%  '
Set-TextValue $ws.Range('D44') '1.843.23'
Set-TextValue $ws.Range('E44') '  +1.51%  '
Set-TextValue $ws.Range('D45') '0.00000000109'
Set-TextValue $ws.Range('E45') '  -2.59%  '
Set-TextValue $ws.Range('D46') '57.41'
Set-TextValue $ws.Range('E46') '  +1.86%  '
Set-TextValue $ws.Range('D47') '8.226'
Set-TextValue $ws.Range('E47') '  +1.10%  '
Set-TextValue $ws.Range('D48') '1.007'
Set-TextValue $ws.Range('E48') '  +0.36%  '
Set-TextValue $ws.Range('D49') '0.05280'
Set-TextValue $ws.Range('E49') '  +0.84%  '
Set-TextValue $ws.Range('D50') '0.4310'
Set-TextValue $ws.Range('E50') '  +0.69%  '
Set-TextValue $ws.Range('D51') '6.069'
Set-TextValue $ws.Range('E51') '  +1.21%  '
